$d = $word.ActiveDocument

# --- Field 1 (paragraph 2): "m:let v = self.name" -> literal text "{m:let v = self.name}" ---
$f1 = $d.Fields.Item(1)
$f1.Delete()
$p2 = $d.Paragraphs(2).Range
$p2.InsertAfter("{m:let v = self.name}")

# --- Field 2 (paragraph 3): " m:v " -> literal text "{m:v}" ---
$f2 = $d.Fields.Item(1)
$f2Start = $f2.Code.Start
$f2.Delete()
$insertPoint = $d.Range($f2Start - 1, $f2Start - 1)
$insertPoint.InsertAfter("{m:v}")

# --- Remove the trailing 4-space run at the end of the last paragraph ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count).Range
$trailSpaces = $d.Range($lastPara.End - 5, $lastPara.End - 1)
$trailSpaces.Delete()
